# Update odds values in row 2 and row 5 of the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("G2").Value = 3.2
$ws.Range("I2").Value = 2.45
$ws.Range("J2").Value = 3.75
$ws.Range("L2").Value = 3.25
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2.05
$ws.Range("S2").Value = 2.4
$ws.Range("T2").Value = 1.53
$ws.Range("W2").Value = 4.5
$ws.Range("X2").Value = 1.18
$ws.Range("AA2").Value = 1.95
$ws.Range("AB2").Value = 1.8
$ws.Range("AC2").Value = 8
$ws.Range("AN2").Value = 7
$ws.Range("AO2").Value = 11

# Row 5 changes
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 5
$ws.Range("N5").Value = 9.5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("AI5").Value = 9.5
$ws.Range("AJ5").Value = 7.5
$ws.Range("AN5").Value = 12
$ws.Range("AO5").Value = 26
$ws.Range("AP5").Value = 17
